$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 518.9032
$ws.Range("I15").Value = 518.9032
$ws.Range("K15").Value = 1556.7096
$ws.Range("M15").Value = -1387.7096

$ws.Range("H33").Value = 174.5
$ws.Range("I33").Value = 170.8
$ws.Range("J33").Value = 230
$ws.Range("K33").Value = 170.8
$ws.Range("L33").Value = 230
$ws.Range("M33").Value = 58.19999999999999
$ws.Range("N33").Value = -688

$ws.Range("H93").Value = 25400.666
$ws.Range("J93").Value = 25400.666
$ws.Range("L93").Value = 25400.666
$ws.Range("N93").Value = -30392.666

$ws.Range("H137").Value = 4357.2104
$ws.Range("I137").Value = 3968.9
$ws.Range("K137").Value = 11906.7
$ws.Range("M137").Value = -9356.700000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3794.0513
$ws.Range("I32").Value = 1332.5
$ws.Range("K32").Value = 1332.5
$ws.Range("M32").Value = -1045.5

$ws.Range("H45").Value = 723.25
$ws.Range("I45").Value = 297.66666
$ws.Range("K45").Value = 297.66666
$ws.Range("M45").Value = 79.33334000000002

$ws.Range("H97").Value = 617.5
$ws.Range("I97").Value = 617.5
$ws.Range("K97").Value = 617.5
$ws.Range("M97").Value = -121.5

$ws.Range("H132").Value = 2336.4443
$ws.Range("I132").Value = 2316
$ws.Range("K132").Value = 6948
$ws.Range("M132").Value = -4418

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 152.16667
$ws.Range("I80").Value = 76.333336
$ws.Range("K80").Value = 76.333336
$ws.Range("M80").Value = 921.666664

$ws.Range("H82").Value = 18521.334
$ws.Range("I82").Value = 14695.111
$ws.Range("K82").Value = 14695.111
$ws.Range("M82").Value = -14312.111

$ws.Range("H83").Value = 152.16667
$ws.Range("I83").Value = 76.333336
$ws.Range("K83").Value = 381.66668
$ws.Range("M83").Value = 4610.33332

$ws.Range("H85").Value = 18521.334
$ws.Range("I85").Value = 14695.111
$ws.Range("K85").Value = 14695.111
$ws.Range("M85").Value = -13369.111

$ws.Range("H94").Value = 1956.4667
$ws.Range("I94").Value = 2102.8462
$ws.Range("J94").Value = 1005
$ws.Range("K94").Value = 2102.8462
$ws.Range("L94").Value = 1005
$ws.Range("M94").Value = -1651.8462
$ws.Range("N94").Value = -1907

$ws.Range("H97").Value = 14212.5
$ws.Range("I97").Value = 14212.5
$ws.Range("K97").Value = 14212.5
$ws.Range("M97").Value = -13221.5

$ws.Range("H105").Value = 1987.8125
$ws.Range("I105").Value = 1961
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1961
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -214
$ws.Range("N105").Value = -5494

$ws.Range("H107").Value = 575.06665
$ws.Range("I107").Value = 485.5
$ws.Range("J107").Value = 933.3333
$ws.Range("K107").Value = 485.5
$ws.Range("L107").Value = 933.3333
$ws.Range("M107").Value = 1434.5
$ws.Range("N107").Value = -4773.3333

$ws.Range("H134").Value = 9352.5
$ws.Range("I134").Value = 8545.857
$ws.Range("K134").Value = 25637.571
$ws.Range("M134").Value = -23102.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2293.9565
$ws.Range("I31").Value = 1867.5
$ws.Range("J31").Value = 2759.182
$ws.Range("K31").Value = 1867.5
$ws.Range("L31").Value = 2759.182
$ws.Range("M31").Value = -1572.5
$ws.Range("N31").Value = -3349.182

$ws.Range("H34").Value = 2293.9565
$ws.Range("I34").Value = 1867.5
$ws.Range("J34").Value = 2759.182
$ws.Range("K34").Value = 1867.5
$ws.Range("L34").Value = 2759.182
$ws.Range("M34").Value = -1665.5
$ws.Range("N34").Value = -3163.182

$ws.Range("H132").Value = 2452.4
$ws.Range("I132").Value = 420.66666
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 1261.99998
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = 1268.00002
$ws.Range("N132").Value = -21560

$ws.Range("H134").Value = 7126
$ws.Range("I134").Value = 6407.5
$ws.Range("K134").Value = 19222.5
$ws.Range("M134").Value = -16687.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1045.125
$ws.Range("J34").Value = 2350
$ws.Range("L34").Value = 7050
$ws.Range("N34").Value = -7218

$ws.Range("H68").Value = 899
$ws.Range("I68").Value = 866
$ws.Range("K68").Value = 2598
$ws.Range("M68").Value = -1787

$ws.Range("H71").Value = 899
$ws.Range("I71").Value = 866
$ws.Range("K71").Value = 7794
$ws.Range("M71").Value = -3738

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 83339650
$ws.Range("I70").Value = 83339650
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 83339650
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -83339380
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 83339650
$ws.Range("I73").Value = 83339650
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 83339650
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -83338714
$ws.Range("N73").ClearContents()

$ws.Range("H92").Value = 2375.5
$ws.Range("J92").Value = 2375.5
$ws.Range("L92").Value = 2375.5
$ws.Range("N92").Value = -6119.5

$ws.Range("H102").Value = 2430.5
$ws.Range("I102").Value = 2279.9412
$ws.Range("J102").Value = 4990
$ws.Range("K102").Value = 2279.9412
$ws.Range("L102").Value = 4990
$ws.Range("M102").Value = -657.9412000000002
$ws.Range("N102").Value = -8234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 7785.2
$ws.Range("I32").Value = 7785.2
$ws.Range("K32").Value = 7785.2
$ws.Range("M32").Value = -7468.2

$ws.Range("H40").Value = 4991.2
$ws.Range("I40").Value = 4150.3335
$ws.Range("K40").Value = 4150.3335
$ws.Range("M40").Value = -4014.3335

$ws.Range("H61").Value = 3876.7778
$ws.Range("I61").Value = 3556
$ws.Range("K61").Value = 3556
$ws.Range("M61").Value = -3354

$ws.Range("H93").Value = 2700
$ws.Range("I93").Value = 3266.6667
$ws.Range("K93").Value = 3266.6667
$ws.Range("M93").Value = -2018.6667

$ws.Range("H113").Value = 3876.7778
$ws.Range("I113").Value = 3556
$ws.Range("K113").Value = 3556
$ws.Range("M113").Value = -1386

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 4057.6365
$ws.Range("I136").Value = 2957
$ws.Range("K136").Value = 8871
$ws.Range("M136").Value = -6321

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 998.6667
$ws.Range("I81").Value = 998.6667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1997.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -936.3334
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 998.6667
$ws.Range("I84").Value = 998.6667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9986.666999999999
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4682.666999999999
$ws.Range("N84").ClearContents()

$ws.Range("H132").Value = 3499
$ws.Range("I132").Value = 1999.5
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 5998.5
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -3468.5
$ws.Range("N132").Value = -20055.5
